# Removing any reference on tweets or epitweetr:
# replace the word "tweet" with "post" in every text cell of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $v = $cell.Value2
    if ($v -ne $null -and $v -is [string] -and $v.Contains("tweet")) {
        $cell.Value2 = $v.Replace("tweet", "post")
    }
}

# Match the author's resulting selection (active cell M2)
$ws.Range("M2").Select()
